$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.240.43"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "3.595.37"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "207.78"
$ws.Range("E5").Value = "  +7.81%  "
$ws.Range("D6").Value = "569.77"
$ws.Range("E6").Value = "  -1.70%  "
$ws.Range("D7").Value = "0.613"
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("D9").Value = "0.684"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "63.48"
$ws.Range("E10").Value = "  +12.76%  "
$ws.Range("D11").Value = "0.148"
$ws.Range("E11").Value = "  -2.21%  "
$ws.Range("D12").Value = "0.0000282"
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").Value = "10.37"
$ws.Range("E13").Value = "  +4.64%  "
$ws.Range("D14").Value = "4.162.87"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "3.587.57"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "19.26"
$ws.Range("E16").Value = "  +4.35%  "
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Value = "68.016.86"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "12.24"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").Value = "1.07"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "404.02"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "4.17"
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").Value = "12.39"
$ws.Range("E23").Value = "  +8.32%  "
$ws.Range("D24").Value = "84.82"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").Value = "12.52"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "3.88"
$ws.Range("E27").Value = "  +6.09%  "
$ws.Range("E28").Value = "  +3.17%  "
$ws.Range("D29").Value = "7.61"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.60"
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").Value = "693.29"
$ws.Range("E31").Value = "  +9.12%  "
$ws.Range("D32").Value = "12.17"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").Value = "63.52"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").Value = "41.65"
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("D36").Value = "0.411"
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "3.25"
$ws.Range("E38").Value = "  +8.47%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0760"
$ws.Range("E39").Value = "  -5.05%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "3.16"
$ws.Range("E40").Value = "  +18.18%  "
$ws.Range("D41").Value = "3.161.62"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("D43").Value = "0.997"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "2.68"
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("D45").Value = "2.76"
$ws.Range("E45").Value = "  +8.48%  "
$ws.Range("D46").Value = "0.0413"
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").Value = "3.14"
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("D48").Value = "0.131"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "8.73"
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").Value = "138.95"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").Value = "2.71"
$ws.Range("E51").Value = "  -1.81%  "
